# Resume and Skill Update
#
# Inserts a new bulleted skills line ("Version Control & Repository: Git &
# Github") immediately before the existing "MS Office (MS Word, MS Excel,
# MS Powerpoint)" bullet, and relocates the (hidden) "_GoBack" bookmark
# — which used to sit in the blank paragraph right after "MS Office (...)"
# — onto the end of the newly inserted line, matching where Word leaves it
# after the most recent edit.

$d = $word.ActiveDocument

# 1. Find the "MS Office (...)" bullet paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("MS Office (")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'MS Office (...)' paragraph"
}

$msOfficePara = $d.Paragraphs.Item($targetIndex)

# 2. Insert a new paragraph right before it. Word copies the paragraph's
#    pPr/rPr (pStyle ListParagraph, numPr bullet, bold rPr) onto the new,
#    still-empty paragraph, so it keeps the same bullet formatting.
$msOfficePara.Range.InsertParagraphBefore()

# 3. Fill in the text of the new paragraph (it now lives at $targetIndex,
#    "MS Office (...)" shifted down by one).
$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "Version Control & Repository: Git & Github"

# 4. Move the "_GoBack" bookmark from the trailing blank paragraph onto the
#    end of the new line (i.e. right before its paragraph mark), which is
#    where Word leaves "_GoBack" after inserting/editing text there.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newPara2 = $d.Paragraphs.Item($targetIndex)
$eopPos = $newPara2.Range.End - 1

# Adding a bookmark with a collapsed range exactly at a paragraph boundary
# is unreliable, so briefly insert a one-character placeholder, bookmark
# the (now non-boundary) range around it, then delete the placeholder —
# the bookmark collapses cleanly to the intended position.
$d.Range($eopPos, $eopPos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($eopPos, $eopPos + 1))
$d.Range($eopPos, $eopPos + 1).Text = ""
